$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column H values for rows 873-898 (no F/G changes)
$ws.Cells.Item(873, 8).Value = 723
$ws.Cells.Item(874, 8).Value = 744
$ws.Cells.Item(875, 8).Value = 761
$ws.Cells.Item(876, 8).Value = 740
$ws.Cells.Item(877, 8).Value = 716
$ws.Cells.Item(878, 8).Value = 724
$ws.Cells.Item(879, 8).Value = 752
$ws.Cells.Item(880, 8).Value = 723
$ws.Cells.Item(881, 8).Value = 700
$ws.Cells.Item(882, 8).Value = 656
$ws.Cells.Item(883, 8).Value = 660
$ws.Cells.Item(884, 8).Value = 587
$ws.Cells.Item(885, 8).Value = 613
$ws.Cells.Item(886, 8).Value = 651
$ws.Cells.Item(887, 8).Value = 600
$ws.Cells.Item(888, 8).Value = 563
$ws.Cells.Item(889, 8).Value = 548
$ws.Cells.Item(890, 8).Value = 529
$ws.Cells.Item(891, 8).Value = 490
$ws.Cells.Item(892, 8).Value = 514
$ws.Cells.Item(893, 8).Value = 540
$ws.Cells.Item(894, 8).Value = 514
$ws.Cells.Item(895, 8).Value = 509
$ws.Cells.Item(896, 8).Value = 509
$ws.Cells.Item(897, 8).Value = 521
$ws.Cells.Item(898, 8).Value = 495

# Update F/G/H values for rows 899-921
$ws.Cells.Item(899, 6).Value = 922
$ws.Cells.Item(899, 8).Value = 517
$ws.Cells.Item(900, 6).Value = 1217
$ws.Cells.Item(900, 8).Value = 551
$ws.Cells.Item(901, 6).Value = 5068
$ws.Cells.Item(901, 7).Value = 460
$ws.Cells.Item(901, 8).Value = 549
$ws.Cells.Item(902, 6).Value = 3021
$ws.Cells.Item(902, 8).Value = 511
$ws.Cells.Item(903, 6).Value = 3237
$ws.Cells.Item(903, 8).Value = 467
$ws.Cells.Item(904, 8).Value = 472
$ws.Cells.Item(905, 8).Value = 434
$ws.Cells.Item(906, 6).Value = 739
$ws.Cells.Item(906, 8).Value = 436
$ws.Cells.Item(907, 6).Value = 751
$ws.Cells.Item(907, 7).Value = 93
$ws.Cells.Item(907, 8).Value = 477
$ws.Cells.Item(908, 6).Value = 1366
$ws.Cells.Item(908, 8).Value = 494
$ws.Cells.Item(909, 6).Value = 4303
$ws.Cells.Item(909, 8).Value = 447
$ws.Cells.Item(910, 6).Value = 2737
$ws.Cells.Item(910, 7).Value = 248
$ws.Cells.Item(910, 8).Value = 402
$ws.Cells.Item(911, 6).Value = 1221
$ws.Cells.Item(911, 7).Value = 91
$ws.Cells.Item(911, 8).Value = 411
$ws.Cells.Item(912, 6).Value = 2797
$ws.Cells.Item(912, 7).Value = 256
$ws.Cells.Item(912, 8).Value = 370
$ws.Cells.Item(913, 6).Value = 718
$ws.Cells.Item(913, 7).Value = 56
$ws.Cells.Item(913, 8).Value = 389
$ws.Cells.Item(914, 6).Value = 916
$ws.Cells.Item(914, 7).Value = 71
$ws.Cells.Item(914, 8).Value = 407
$ws.Cells.Item(915, 6).Value = 4283
$ws.Cells.Item(915, 7).Value = 396
$ws.Cells.Item(915, 8).Value = 389
$ws.Cells.Item(916, 6).Value = 3014
$ws.Cells.Item(916, 7).Value = 220
$ws.Cells.Item(916, 8).Value = 378
$ws.Cells.Item(917, 6).Value = 2474
$ws.Cells.Item(917, 7).Value = 189
$ws.Cells.Item(917, 8).Value = 369
$ws.Cells.Item(918, 6).Value = 3492
$ws.Cells.Item(918, 7).Value = 238
$ws.Cells.Item(918, 8).Value = 383
$ws.Cells.Item(919, 6).Value = 3828
$ws.Cells.Item(919, 7).Value = 294
$ws.Cells.Item(919, 8).Value = 355
$ws.Cells.Item(920, 6).Value = 1984
$ws.Cells.Item(920, 7).Value = 72
$ws.Cells.Item(920, 8).Value = 366
$ws.Cells.Item(921, 6).Value = 966
$ws.Cells.Item(921, 7).Value = 92
$ws.Cells.Item(921, 8).Value = 379

# Add new F/G/H values for rows 922-927 (previously empty)
$ws.Cells.Item(922, 6).Value = 5833
$ws.Cells.Item(922, 7).Value = 363
$ws.Cells.Item(922, 8).Value = 361
$ws.Cells.Item(923, 6).Value = 2773
$ws.Cells.Item(923, 7).Value = 231
$ws.Cells.Item(923, 8).Value = 361
$ws.Cells.Item(924, 6).Value = 2509
$ws.Cells.Item(924, 7).Value = 241
$ws.Cells.Item(924, 8).Value = 362
$ws.Cells.Item(925, 6).Value = 1040
$ws.Cells.Item(925, 7).Value = 78
$ws.Cells.Item(925, 8).Value = 359
$ws.Cells.Item(926, 6).Value = 2140
$ws.Cells.Item(926, 7).Value = 263
$ws.Cells.Item(926, 8).Value = 359
$ws.Cells.Item(927, 6).Value = 409
$ws.Cells.Item(927, 7).Value = 57
$ws.Cells.Item(927, 8).Value = 359

# Add new row 928 with full data, matching date format style of column A
$ws.Cells.Item(928, 1).Value = 44822
$ws.Cells.Item(928, 2).Value = 1839213
$ws.Cells.Item(928, 3).Value = 248
$ws.Cells.Item(928, 4).Value = 57
$ws.Cells.Item(928, 5).Value = 20436
$ws.Cells.Item(928, 6).Value = 252
$ws.Cells.Item(928, 7).Value = 52
$ws.Cells.Item(928, 8).Value = 359
$ws.Cells.Item(928, 1).NumberFormat = $ws.Cells.Item(927, 1).NumberFormat
